$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '54.357.71'
$ws.Range("E2").Value = '  +0.19%  '

# Row 3
$ws.Range("D3").Value = '2.280.18'
$ws.Range("E3").Value = '  +0.83%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.37%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.85%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.22%  '

# Row 8
$ws.Range("E8").Value = '  +0.46%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0956'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.47%  '

# Row 10
$ws.Range("E10").Value = '  +1.55%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.333'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.46%  '

# Row 12
$ws.Range("E12").Value = '  +1.55%  '

# Row 13
$ws.Range("D13").Value = '2.686.33'
$ws.Range("E13").Value = '  +1.00%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.47%  '

# Row 15
$ws.Range("D15").Value = '54.296.96'
$ws.Range("E15").Value = '  +0.23%  '

# Row 16
$ws.Range("E16").Value = '  +0.46%  '

# Row 17
$ws.Range("D17").Value = '2.302.31'
$ws.Range("E17").Value = '  +1.96%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.28%  '

# Row 19
$ws.Range("E19").Value = '  +2.35%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.00%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '62.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.87%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '

# Row 25
$ws.Range("D25").Value = '2.389.39'
$ws.Range("E25").Value = '  +1.21%  '

# Row 26
$ws.Range("E26").Value = '  +3.19%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '174.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.80%  '

# Row 29
$ws.Range("E29").Value = '  +1.37%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.60%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0686'
$ws.Range("E31").Value = '  +1.04%  '

# Row 32
$ws.Range("E32").Value = '  +1.55%  '

# Row 33
$ws.Range("E33").Value = '  -0.02%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.03%  '

# Row 35
$ws.Range("E35").Value = '  -0.13%  '

# Row 36
$ws.Range("E36").Value = '  +9.64%  '

# Row 37
$ws.Range("E37").Value = '  +0.90%  '

# Row 38
$ws.Range("E38").Value = '  +3.65%  '

# Row 39
$ws.Range("E39").Value = '  +0.06%  '

# Row 40
$ws.Range("E40").Value = '  +1.43%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '126.01'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.63%  '

# Row 43
$ws.Range("E43").Value = '  -0.33%  '

# Row 44
$ws.Range("B44").Value = 'Hedera'
$ws.Range("C44").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0491'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.60%  '

# Row 45
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0897'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.71%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.546'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.81%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '240.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.87%  '

# Row 48
$ws.Range("E48").Value = '  -0.06%  '

# Row 49
$ws.Range("E49").Value = '  +1.51%  '

# Row 50
$ws.Range("E50").Value = '  +0.96%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.83%  '
